$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Translate Publisher contact text from German to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Add Description value (was previously blank)
$ws.Range("B12").Value = "consent states - subset WITHDRAWAL documents"
